# SummonTable.xlsx fix: add "SummonGrade : Int" column in front of the
# existing ItemId/Probability table, duplicate the old 6 rows into a
# second grade-2 block, and add per-block SUM checks in column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Grow the table first (A2:B8 -> A2:C14) --------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A2:C14"))

# --- 2. Header row (also renames the ListColumns) ------------------------
$ws.Range("A2").Value = "SummonGrade : Int"
$ws.Range("B2").Value = "ItemId : Int"
$ws.Range("C2").Value = "Probability : Int"

# --- 3. Grade-1 block (rows 3-8) -----------------------------------------
$ws.Range("A3").Value = 1

$ws.Range("B3").Value = 10001
$ws.Range("B4").Value = 10002
$ws.Range("B5").Value = 10003
$ws.Range("B6").Value = 10005
$ws.Range("B7").Value = 10006
$ws.Range("B8").Value = 10007

$ws.Range("C3").Value = 5000
$ws.Range("C4").Value = 2500
$ws.Range("C5").Value = 1500
$ws.Range("C6").Value = 700
$ws.Range("C7").Value = 250
$ws.Range("C8").Value = 50

# A4 is a lone (non-shared) formula, A5:A8 become one shared formula group
$ws.Range("A4").Formula = "=A3"
$ws.Range("A5:A8").Formula = "=A4"

# Grade-1 total
$ws.Range("D8").Formula = "=SUM(C3:C8)"

# --- 4. Grade-2 block (rows 9-14), mirrors rows 3-8 -----------------------
$ws.Range("A9").Value = 2

$ws.Range("C9").Value = 3500
$ws.Range("C10").Value = 3000
$ws.Range("C11").Value = 2000
$ws.Range("C12").Value = 1000
$ws.Range("C13").Value = 400
$ws.Range("C14").Value = 100

# B9 is a lone formula referencing B3; B10:B14 is a shared formula group
$ws.Range("B9").Formula = "=B3"
$ws.Range("B10:B14").Formula = "=B4"

# A10 is a lone formula referencing A9; A11:A14 is a shared formula group
$ws.Range("A10").Formula = "=A9"
$ws.Range("A11:A14").Formula = "=A10"

# Grade-2 total
$ws.Range("D14").Formula = "=SUM(C9:C14)"

# --- 5. Column widths (closest representable values) ---------------------
$ws.Columns.Item(1).ColumnWidth = 21.142857142857146
$ws.Columns.Item(2).ColumnWidth = 16.857142857142854
$ws.Columns.Item(3).ColumnWidth = 19.142857142857146

# --- 6. Selection matches the post-edit cursor position -------------------
$ws.Range("D10").Select() | Out-Null
